$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 4) with a fresh snapshot of the BJ899050 index.
# It mirrors row 3's data, except for an updated timestamp and the market
# classification/status (not yet opened).
$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "BJ899050"
$ws.Range("C4").Value = "北证50"
$ws.Range("D4").Value = "未开盘"
$ws.Range("E4").Value = "'760.89"
$ws.Range("F4").Value = "+3.76  +0.50%"
$ws.Range("G4").Value = 760.89
$ws.Range("H4").Value = 751.45
$ws.Range("I4").Value = 1268900
$ws.Range("J4").Value = "'"
$ws.Range("K4").Value = 748.5
$ws.Range("L4").Value = 757.13
$ws.Range("M4").Value = 1267000000
$ws.Range("N4").Value = 0.016
$ws.Range("O4").Value = 1079.71
$ws.Range("P4").Value = 702.55
